$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cell E1, matching style/formatting of D1 (bold, border, centered)
$ws.Range("D1").Copy()
$ws.Range("E1").PasteSpecial(-4122)
$ws.Range("E1").Value = 3

# Data values for E2:E22
$values = @(
    -0.001884944920440637,
    0.002281926303347229,
    -0.0004456408577852223,
    0.002325581395348837,
    0.002308815810710199,
    0.002316777923169888,
    -0.001942167266456478,
    0.0008844103768215862,
    -0.001206072474370804,
    0.001031927034490668,
    0.00002214839424141677,
    0.0002972547648190249,
    -0.00002768549280177161,
    -0.0004228329809725159,
    -0.000225616105518917,
    -0.0002474022761009403,
    0.002325581395348837,
    -0.001223990208078335,
    0.0006536505886834703,
    0.0008347688774045363,
    0.002325581395348837
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 5).Value = $values[$i]
}
